# render site on laptop
#
# 1) Collapse the run-per-word title paragraph into a single run with the
#    full text "Modern Dive Chapter 2: Data Visualization".
# 2) Subtitle style: base it on "Title" instead of "Normal", and drop the
#    explicit font color override (-> Automatic).
# 3) AbstractTitle style: drop its explicit font color override (-> Automatic).

$d = $word.ActiveDocument

# --- 1. Merge the Title paragraph's runs into one run -----------------
$titlePara = $d.Paragraphs(1)
$titleText = "Modern Dive Chapter 2: Data Visualization"

$fullRange = $titlePara.Range
$fullRange.MoveEnd(1, -1) | Out-Null        # exclude the paragraph mark
$paraStart = $fullRange.Start
$paraEnd = $fullRange.End

# Delete everything after the first word, then overwrite the first word's
# text with the complete sentence -- this keeps a single <w:r> (the first
# run's own formatting, i.e. none) instead of leaving one run per word.
$firstWordLen = "Modern".Length
if ($paraEnd -gt ($paraStart + $firstWordLen)) {
    $rest = $d.Range($paraStart + $firstWordLen, $paraEnd)
    $rest.Delete()
}
$firstRun = $d.Range($paraStart, $paraStart + $firstWordLen)
$firstRun.Text = $titleText

# --- 2. Subtitle style: re-base on Title, clear direct font color -----
$titleStyle = $d.Styles("Title")
$subtitleStyle = $d.Styles("Subtitle")
$subtitleStyle.BaseStyle = $titleStyle
$subtitleStyle.Font.Color = -16777216   # wdColorAutomatic

# --- 3. AbstractTitle style: clear direct font color -------------------
$abstractTitleStyle = $d.Styles("AbstractTitle")
$abstractTitleStyle.Font.Color = -16777216   # wdColorAutomatic
